$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Populate new cells in the order that reproduces the original
# shared-string insertion order (D4, A5, A4, D5, D6).
$ws.Range("D4").Value = "Change thread count setting"
$ws.Range("A5").Value = "1e8b92d"
$ws.Range("A4").Value = "a88d2af"
$ws.Range("D5").Value = "Make Slice DataDictionary lazy eval"
$ws.Range("D6").Value = "Make algo run on dedicated thread"

$ws.Range("B4").Value = 19.16
$ws.Range("C4").Value = 186

$ws.Range("B5").Value = 18.37
$ws.Range("C5").Value = 195

$ws.Range("B6").Value = 18.3
$ws.Range("C6").Value = 196

$ws.Range("D6").Select()
